$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 4630246
$ws.Range("I19").Value = 7937280.5
$ws.Range("J19").Value = 397.8
$ws.Range("K19").Value = 7937280.5
$ws.Range("L19").Value = 397.8
$ws.Range("M19").Value = -7937105.5
$ws.Range("N19").Value = -747.8

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H52").Value = 519.625
$ws.Range("I52").Value = 219
$ws.Range("K52").Value = 657
$ws.Range("M52").Value = -497

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 2367.818
$ws.Range("I132").Value = 2290.7856
$ws.Range("J132").Value = 2799.2
$ws.Range("K132").Value = 6872.3568
$ws.Range("L132").Value = 8397.599999999999
$ws.Range("M132").Value = -4342.3568
$ws.Range("N132").Value = -13457.6

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 5856962.5
$ws.Range("I138").Value = 11891.5
$ws.Range("J138").Value = 12351486
$ws.Range("K138").Value = 35674.5
$ws.Range("L138").Value = 37054458
$ws.Range("M138").Value = -30534.5
$ws.Range("N138").Value = -37064738

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 3780.9333
$ws.Range("I2").Value = 2739.25
$ws.Range("J2").Value = 4971.4287
$ws.Range("K2").Value = 2739.25
$ws.Range("L2").Value = 4971.4287
$ws.Range("M2").Value = -2626.25
$ws.Range("N2").Value = -5197.4287

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H34").Value = 10512.5
$ws.Range("J34").Value = 0
$ws.Range("L34").Value = 0
$ws.Range("N34").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H36").Value = 4034
$ws.Range("I36").Value = 4034
$ws.Range("K36").Value = 4034
$ws.Range("M36").Value = -3688

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H51").Value = 25000
$ws.Range("J51").Value = 25000
$ws.Range("L51").Value = 25000
$ws.Range("N51").Value = -26512

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 6423.619
$ws.Range("I61").Value = 6205.1055
$ws.Range("J61").Value = 8499.5
$ws.Range("K61").Value = 6205.1055
$ws.Range("L61").Value = 8499.5
$ws.Range("M61").Value = -5993.1055
$ws.Range("N61").Value = -8923.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H116").Value = 3780.9333
$ws.Range("I116").Value = 2739.25
$ws.Range("J116").Value = 4971.4287
$ws.Range("K116").Value = 2739.25
$ws.Range("L116").Value = 4971.4287
$ws.Range("M116").Value = -445.25
$ws.Range("N116").Value = -9559.4287

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 6423.619
$ws.Range("I136").Value = 6205.1055
$ws.Range("J136").Value = 8499.5
$ws.Range("K136").Value = 18615.3165
$ws.Range("L136").Value = 25498.5
$ws.Range("M136").Value = -16065.3165
$ws.Range("N136").Value = -30598.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 3780.9333
$ws.Range("I3").Value = 2739.25
$ws.Range("J3").Value = 4971.4287
$ws.Range("K3").Value = 2739.25
$ws.Range("L3").Value = 4971.4287
$ws.Range("M3").Value = -2625.25
$ws.Range("N3").Value = -5199.4287

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 4196.2383
$ws.Range("I107").Value = 2417.5715
$ws.Range("K107").Value = 2417.5715
$ws.Range("M107").Value = -497.5715

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H132").Value = 158715.38
$ws.Range("J132").Value = 158715.38
$ws.Range("L132").Value = 158715.38
$ws.Range("N132").Value = -168835.38

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H4").Value = 142.5
$ws.Range("I4").Value = 171.5
$ws.Range("J4").Value = 41
$ws.Range("K4").Value = 171.5
$ws.Range("L4").Value = 41
$ws.Range("M4").Value = -59.5
$ws.Range("N4").Value = -265

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H63").Value = 40000
$ws.Range("J63").Value = 40000
$ws.Range("L63").Value = 40000
$ws.Range("N63").Value = -41372

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H66").Value = 40000
$ws.Range("J66").Value = 40000
$ws.Range("L66").Value = 120000
$ws.Range("N66").Value = -126864

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H82").Value = 115000
$ws.Range("J82").Value = 115000
$ws.Range("L82").Value = 115000
$ws.Range("N82").Value = -115722

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H85").Value = 115000
$ws.Range("J85").Value = 115000
$ws.Range("L85").Value = 115000
$ws.Range("N85").Value = -117496

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H103").Value = 14420.667
$ws.Range("I103").Value = 14420.667
$ws.Range("K103").Value = 14420.667
$ws.Range("M103").Value = -13248.667

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 2041.1666
$ws.Range("I132").Value = 2029.5
$ws.Range("K132").Value = 6088.5
$ws.Range("M132").Value = -3558.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H14").Value = 650.8182
$ws.Range("I14").Value = 650.8182
$ws.Range("K14").Value = 1952.4546
$ws.Range("M14").Value = -1779.4546

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H99").Value = 34500
$ws.Range("I99").Value = 35000
$ws.Range("J99").Value = 34000
$ws.Range("K99").Value = 105000
$ws.Range("L99").Value = 102000
$ws.Range("M99").Value = -102754
$ws.Range("N99").Value = -106492

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H107").Value = 2466.111
$ws.Range("I107").Value = 3979.3333
$ws.Range("J107").Value = 1709.5
$ws.Range("K107").Value = 11937.9999
$ws.Range("L107").Value = 5128.5
$ws.Range("M107").Value = -10017.9999
$ws.Range("N107").Value = -8968.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H52").Value = 29010.5
$ws.Range("J52").Value = 33500.75
$ws.Range("L52").Value = 33500.75
$ws.Range("N52").Value = -34018.75

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H117").Value = 44000
$ws.Range("J117").Value = 44000
$ws.Range("L117").Value = 44000
$ws.Range("N117").Value = -50884

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 6110.1113
$ws.Range("I122").Value = 4831.8335
$ws.Range("J122").Value = 8666.666999999999
$ws.Range("K122").Value = 14495.5005
$ws.Range("L122").Value = 26000.001
$ws.Range("M122").Value = -12045.5005
$ws.Range("N122").Value = -30900.001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H2").Value = 200
$ws.Range("I2").Value = 200
$ws.Range("K2").Value = 200
$ws.Range("M2").Value = -88

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 5536.6665
$ws.Range("I40").Value = 2333.3333
$ws.Range("K40").Value = 2333.3333
$ws.Range("M40").Value = -2197.3333

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H75").Value = 400021600
$ws.Range("J75").Value = 400021600
$ws.Range("L75").Value = 400021600
$ws.Range("N75").Value = -400023472

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H78").Value = 400021600
$ws.Range("J78").Value = 400021600
$ws.Range("L78").Value = 1200064800
$ws.Range("N78").Value = -1200074160

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 4156.357
$ws.Range("I122").Value = 3754.4443
$ws.Range("J122").Value = 4879.8
$ws.Range("K122").Value = 11263.3329
$ws.Range("L122").Value = 14639.4
$ws.Range("M122").Value = -8813.332900000001
$ws.Range("N122").Value = -19539.4

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H129").Value = 205949.5
$ws.Range("J129").Value = 205949.5
$ws.Range("L129").Value = 205949.5
$ws.Range("N129").Value = -215949.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 4029.8333
$ws.Range("I132").Value = 3643
$ws.Range("K132").Value = 10929
$ws.Range("M132").Value = -8399

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 15467.262
$ws.Range("I136").Value = 3770.84
$ws.Range("J136").Value = 32667.883
$ws.Range("K136").Value = 11312.52
$ws.Range("L136").Value = 98003.649
$ws.Range("M136").Value = -8762.52
$ws.Range("N136").Value = -103103.649

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H14").Value = 24
$ws.Range("J14").Value = 24
$ws.Range("L14").Value = 24
$ws.Range("N14").Value = -360

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 3201.772
$ws.Range("I122").Value = 1990.75
$ws.Range("K122").Value = 5972.25
$ws.Range("M122").Value = -3522.25

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H129").Value = 68200
$ws.Range("J129").Value = 68200
$ws.Range("L129").Value = 68200
$ws.Range("N129").Value = -78200
